# Updated to 1/25/2020 for Wuhan nCoV
#
# The "数据" (data) sheet used to carry an SIR-style analysis (columns
# I / dI/dt / 新增病例-dI/dt / dR/dt in J:M) alongside the raw daily
# case counts in A:I. This update drops that analysis block entirely,
# keeping only the raw daily figures in columns A:I (J18 is left as an
# empty, still-styled cell, matching the author's original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the now-unused analysis columns (J:M) for every data row.
# ClearContents (unlike Delete) does not shift any cells, which is what
# we want here -- the J18 cell keeps its existing style (s="6") as an
# empty cell, and the now-unreferenced shared strings ("I", "dI/dt",
# "新增病例-dI/dt", "dR/dt") fall out of the shared-strings table on save.
$ws.Range("J1:M35").ClearContents()

# Move the saved selection to J18, matching the author's final cursor
# position after trimming the sheet down to the raw-data columns.
[void]$ws.Range("J18").Select()
